$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.119.06'
$ws.Range("E2").Value = '  -1.09%  '
$ws.Range("D3").Value = '2.299.88'
$ws.Range("E3").Value = '  -2.52%  '
$ws.Range("E4").Value = '  -0.02%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '312.41'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -3.48%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '104.52'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +2.69%  '
$ws.Range("E7").Value = '  -1.88%  '
$ws.Range("E8").Value = '  +0.01%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.607'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.87%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '40.31'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.09%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0913'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.98%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '8.28'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -2.12%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.973'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -2.99%  '
$ws.Range("E15").Value = '  -4.81%  '
$ws.Range("D16").Value = '2.648.02'
$ws.Range("E16").Value = '  -2.53%  '
$ws.Range("D17").Value = '2.300.15'
$ws.Range("E17").Value = '  -2.54%  '
$ws.Range("D18").Value = '42.005.47'
$ws.Range("E18").Value = '  -1.67%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '7.61'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -5.65%  '
$ws.Range("E20").Value = '  -1.64%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '74.65'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.91%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '3.46'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -6.20%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '257.79'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -2.93%  '
$ws.Range("E24").Value = '  -0.46%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '9.31'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -8.16%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.32%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '10.98'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -4.28%  '
$ws.Range("E28").Value = '  +3.27%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '22.80'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.70%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '166.26'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -5.43%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '35.75'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.19%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.0895'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.87%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '2.91'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -5.74%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '5.83'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -4.06%  '
$ws.Range("E35").Value = '  +10.84%  '
$ws.Range("E36").Value = '  -2.61%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '4.55'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.44%  '
$ws.Range("E38").Value = '  -1.29%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.75'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -6.14%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '3.62'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -4.18%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '72.04'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +2.82%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '98.51'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +8.14%  '
$ws.Range("E43").Value = '  -3.52%  '
$ws.Range("E44").Value = '  -4.89%  '
$ws.Range("E45").Value = '  -0.07%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '12.35'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +3.58%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '112.05'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -7.23%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '9.04'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -2.22%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '5.33'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -3.86%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '74.83'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +4.78%  '
$ws.Range("D51").Value = '1.570.72'
$ws.Range("E51").Value = '  +0.69%  '
